# [LPF-879]: CCMS Third party report
# Remove the "By Source and Expenditure type", "Provider Contigency" and
# "MAIN" worksheets, leaving only "Summary" and "Transparency Rec".

$wb = $excel.ActiveWorkbook

$excel.DisplayAlerts = $false

# "By Source and Expenditure type" hosts its own pivot table ("SourceNType")
# built on the same pivot cache as the one kept on "Transparency Rec". Clear
# the pivot table explicitly before removing the sheet so the now-unused
# pivot table definition doesn't linger as an orphaned part.
$wsBySource = $wb.Worksheets.Item("By Source and Expenditure type")
$pt = $wsBySource.PivotTables().Item(1)
$pt.TableRange2.Clear() | Out-Null
$wsBySource.Delete() | Out-Null

$wb.Worksheets.Item("Provider Contigency").Delete() | Out-Null
$wb.Worksheets.Item("MAIN").Delete() | Out-Null

$excel.DisplayAlerts = $true
